$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Nom"
$ws.Range("B1").Value = "Kill"
$ws.Range("C1").Value = "Points"

$ws.Range("A2").Value = "Come"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3

$ws.Range("E7").Select()
